$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing row down to the new row
$ws.Range("A75:F75").Copy()
$ws.Range("A76:F76").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's values
$ws.Range("A76").Value = 43978
$ws.Range("B76").Value = 616
$ws.Range("C76").Value = 239
$ws.Range("D76").Value = 467
$ws.Range("E76").Value = 25
$ws.Range("F76").Value = 27

# Resize the table to include the new row
$ws.ListObjects.Item("Condicion_Pacientes").Resize($ws.Range("A1:F76"))

# Update selection to match the edited cell and scroll the window down
# so the new row is visible (mirrors the saved view state)
$ws.Range("F76").Select()
$excel.ActiveWindow.ScrollRow = 67
$excel.ActiveWindow.ScrollColumn = 1
